$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 125
$ws.Range("F5").Value = 577
$ws.Range("F6").Value = 1238
$ws.Range("F7").Value = 1021
$ws.Range("F10").Value = 749
$ws.Range("F13").Value = 864
$ws.Range("F14").Value = 23
$ws.Range("F15").Value = 94
$ws.Range("F16").Value = 1332
$ws.Range("F18").Value = 29
$ws.Range("F20").Value = 102
$ws.Range("F21").Value = 32
$ws.Range("F22").Value = 1253
$ws.Range("F23").Value = 355
$ws.Range("F26").Value = 120
$ws.Range("F27").Value = 2512
$ws.Range("F41").Value = 211

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 51
$ws.Range("F5").Value = 675
$ws.Range("F6").Value = 491
$ws.Range("F7").Value = 5
$ws.Range("F11").Value = 251
$ws.Range("F16").Value = 647
$ws.Range("F18").Value = 19
$ws.Range("F19").Value = 486
$ws.Range("F21").Value = 13

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2030
$ws.Range("F6").Value = 2185
$ws.Range("F7").Value = 859
$ws.Range("F8").Value = 837
$ws.Range("F11").Value = 955
$ws.Range("F12").Value = 170
$ws.Range("F13").Value = 37

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2030
$ws.Range("F4").Value = 2185
$ws.Range("F8").Value = 859
$ws.Range("F9").Value = 837
$ws.Range("F11").Value = 577
$ws.Range("F12").Value = 1238
$ws.Range("F13").Value = 1021
$ws.Range("F14").Value = 955
$ws.Range("F16").Value = 749
$ws.Range("F19").Value = 170
$ws.Range("F20").Value = 37
$ws.Range("F21").Value = 675
$ws.Range("F22").Value = 864
$ws.Range("F23").Value = 23
$ws.Range("F24").Value = 94
$ws.Range("F25").Value = 1332
$ws.Range("F26").Value = 491
$ws.Range("F28").Value = 29
$ws.Range("F30").Value = 102
$ws.Range("F31").Value = 32
$ws.Range("F32").Value = 1253
$ws.Range("F33").Value = 355
$ws.Range("F35").Value = 120
$ws.Range("F36").Value = 2512
$ws.Range("F41").Value = 251
$ws.Range("F46").Value = 19
$ws.Range("F47").Value = 13
$ws.Range("F50").Value = 211
